$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range (B:E) keeps its original text formatting so that
# numeric-looking values (e.g. "569.22") are stored as text, not numbers,
# matching the source workbook which stores every cell as an inline string.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.932.74'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '2.732.79'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '569.22'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').Value = '158.21'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.595'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('E10').Value = '  +4.21%  '
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').Value = '0.381'
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').Value = '3.217.51'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').Value = '26.53'
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').Value = '63.576.34'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('D17').Value = '2.738.28'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').Value = '12.03'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = '4.79'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('D20').Value = '352.95'
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').Value = '0.522'
$ws.Range('E23').Value = '  -6.95%  '
$ws.Range('D24').Value = '64.23'
$ws.Range('E24').Value = '  -2.86%  '
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('D29').Value = '1.96'
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '1.35'
$ws.Range('E30').Value = '  +6.65%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').Value = '162.89'
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '20.01'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '4.88'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').Value = '0.985'
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '349.41'
$ws.Range('E39').Value = '  +4.70%  '
$ws.Range('D40').Value = '6.27'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').Value = '4.11'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('D42').Value = '38.57'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = '21.99'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = '21.13'
$ws.Range('E44').Value = '  -3.82%  '
$ws.Range('D45').Value = '0.0581'
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('D46').Value = '134.76'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('D47').Value = '0.621'
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('D49').Value = '0.0247'
$ws.Range('E49').Value = '  -4.13%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Value = '11.03'
$ws.Range('E51').Value = '  -0.12%  '
